# QUERY_10set_1times.xlsx — drop the "QUERY" label row so the QUERY(...)
# formula (previously in A2) becomes the sheet's only cell, at A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 held the shared-string label "QUERY"; deleting it shifts the
# formula row (old A2) up into A1 and shrinks the used range to just A1.
$ws.Rows.Item(1).Delete()

# Column A widens a bit now that its only content is the long formula
# text (the sheet shows formulas, not values).
$ws.Columns.Item(1).ColumnWidth = 31.8333333333333

# Selection moves to A2, just below the remaining formula cell.
$ws.Range("A2").Select()

# Best-effort: restore the saved window's vertical scroll position.
try { $excel.ActiveWindow.Top = 1905 } catch {}
